$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (O) by copying column N's row-by-row formatting
# (borders/fonts/number formats) and values for rows 3-14, then overwrite
# the cells whose values actually differ for the new year.
$ws.Range("N3:N14").Copy($ws.Range("O3:O14")) | Out-Null

# Header year
$ws.Range("O4").Value2 = 2021

# Per-indicator values for 2021 (rows 5-14); most repeat the 2020 figure,
# a few move.
$ws.Range("O5").Value2 = 97
$ws.Range("O6").Value2 = 96.2
$ws.Range("O7").Value2 = 62.7
$ws.Range("O8").Value2 = 100
$ws.Range("O9").Value2 = 100
$ws.Range("O10").Value2 = "-"
$ws.Range("O11").Value2 = 100
$ws.Range("O12").Value2 = 57.9
$ws.Range("O13").Value2 = 100
$ws.Range("O14").Value2 = "-"

# Match the saved selection from the source edit (moved off the grid below
# the data, column O).
$ws.Range("O17").Select() | Out-Null
